$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New IDs for column A (was l01..l05, becomes LSH001..LSH005)
$ids = @("LSH001", "LSH002", "LSH003", "LSH004", "LSH005")
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $ids[$i]
}

# Remove the old teacher column (C: GV001..GV006) entirely, shifting
# the CNTT column (D) left into column C.
$ws.Range("C1:C5").Delete()

# Update the active selection to match the saved workbook state.
$ws.Range("F6").Select()
